$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 19, shifting the existing
# rows 19-24 down to 21-26.
$ws.Rows("19:20").Insert()

# New row 19: weekly update record (Región Metropolitana, calidad Primera)
$ws.Cells.Item(19, 1).Value = 4
$ws.Cells.Item(19, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(19, 3).Value = "Los Lagos"
$ws.Cells.Item(19, 4).Value = 44582
$ws.Cells.Item(19, 5).Value = 10
$ws.Cells.Item(19, 6).Value = 100112030
$ws.Cells.Item(19, 7).Value = "Poroto granado"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 35000
$ws.Cells.Item(19, 12).Value = 35000
$ws.Cells.Item(19, 13).Value = 35000
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 1400
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# New row 20: weekly update record (Región Metropolitana, calidad Segunda)
$ws.Cells.Item(20, 1).Value = 4
$ws.Cells.Item(20, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(20, 3).Value = "Los Lagos"
$ws.Cells.Item(20, 4).Value = 44582
$ws.Cells.Item(20, 5).Value = 10
$ws.Cells.Item(20, 6).Value = 100112030
$ws.Cells.Item(20, 7).Value = "Poroto granado"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Segunda"
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 27000
$ws.Cells.Item(20, 12).Value = 27000
$ws.Cells.Item(20, 13).Value = 27000
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1080
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# Match the date-number formatting used by the other rows in column D.
$ws.Range("D19:D20").NumberFormat = $ws.Range("D21").NumberFormat
